# Admin beach creating/editing/removing/importing. Also fixed bug in beach import.
# Append four newly-imported beach rows (106-109) to the sample beach list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,6

# Row 38 - BeachID 106
$data[0,0] = 106
$data[0,1] = "wdjiasioj"
$data[0,2] = "dsads"
$data[0,3] = "zfd"
$data[0,4] = 4
$data[0,5] = 2

# Row 39 - BeachID 107
$data[1,0] = 107
$data[1,1] = "EFS"
$data[1,2] = "2DASD"
$data[1,3] = "DSADSA"
$data[1,4] = 3
$data[1,5] = 3

# Row 40 - BeachID 108
$data[2,0] = 108
$data[2,1] = "dsadas"
$data[2,2] = "dsads"
$data[2,3] = "adsas"
$data[2,4] = 2
$data[2,5] = 2342

# Row 41 - BeachID 109
$data[3,0] = 109
$data[3,1] = "feaf"
$data[3,2] = "fds"
$data[3,3] = "fsdfd"
$data[3,4] = 3
$data[3,5] = 4

$ws.Range("A38:F41").Value = $data

# Move the view/selection to the newly added last row, like the source file shows.
$ws.Range("F41").Select()
